$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E is empty (A:D has data, then F has the "Time" column).
# Delete the empty column E so column F's data shifts left into E.
$ws.Range("E1").EntireColumn.Delete()

# Update the active selection to match the saved view state.
$ws.Range("F11").Select()
